$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

$ws.Range("A4").Value = "`${msg.getProperty('search_full_text')}"
$ws.Range("B4").Value = "`${search_full_text}"
